$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 229
$ws.Range("F7").Value = 12815
$ws.Range("G7").Value = 238
$ws.Range("F8").Value = 36
$ws.Range("F9").Value = 114
$ws.Range("F10").Value = 204
$ws.Range("F11").Value = 2879
$ws.Range("F13").Value = 6157
$ws.Range("F16").Value = 3276
$ws.Range("F17").Value = 22
$ws.Range("F18").Value = 155
$ws.Range("F21").Value = 30
$ws.Range("F22").Value = 52
$ws.Range("F24").Value = 3537
$ws.Range("F26").Value = 86
$ws.Range("F27").Value = 2641
$ws.Range("F28").Value = 357
$ws.Range("F29").Value = 1839
$ws.Range("F30").Value = 96
$ws.Range("F32").Value = 6425
$ws.Range("C33").Value = "北京·DICE CON 2024 第八届国际桌面游戏展"
$ws.Range("F33").Value = 12
$ws.Range("F35").Value = 124
$ws.Range("F36").Value = 1933
$ws.Range("F37").Value = 1285
$ws.Range("F38").Value = 77
$ws.Range("F39").Value = 999
$ws.Range("F41").Value = 195
$ws.Range("F42").Value = 212
$ws.Range("F45").Value = 117
$ws.Range("F46").Value = 1162
$ws.Range("F47").Value = 1709
$ws.Range("F48").Value = 143

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 90
$ws.Range("F24").Value = 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 404
$ws.Range("F3").Value = 565

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 404
$ws.Range("F7").Value = 565
$ws.Range("F8").Value = 229
$ws.Range("F11").Value = 12815
$ws.Range("G11").Value = 238
$ws.Range("F12").Value = 114
$ws.Range("F14").Value = 204
$ws.Range("F15").Value = 2879
$ws.Range("F17").Value = 6157
$ws.Range("F20").Value = 30
$ws.Range("F21").Value = 52
$ws.Range("F24").Value = 3537
$ws.Range("F27").Value = 2641
$ws.Range("F28").Value = 1839
$ws.Range("F29").Value = 96
$ws.Range("F31").Value = 6425
$ws.Range("F32").Value = 90
$ws.Range("F34").Value = 124
$ws.Range("F35").Value = 1933
$ws.Range("F37").Value = 1285
$ws.Range("F38").Value = 77
$ws.Range("F39").Value = 999
$ws.Range("F40").Value = 195
$ws.Range("F41").Value = 212
$ws.Range("F44").Value = 117
$ws.Range("F45").Value = 1162
$ws.Range("F47").Value = 1709
$ws.Range("F48").Value = 143
